# Analyzing nodetype by name
# Applies the cell-content / formatting changes described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells that lost their content entirely ---
$ws.Range("G1").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("D12").Value = ""

# --- Update cell text values ---
$ws.Range("B2").Value = "父级（可自定义）"
$ws.Range("D3").Value = "b_"
$ws.Range("D6").Value = "b_"
$ws.Range("F6").Value = "Students@ScrollView:V"
$ws.Range("C7").Value = "l,r,t,b"
$ws.Range("F7").Value = "leftBar@Scrollbar:B"
$ws.Range("F8").Value = "Progress@Slider:B"
$ws.Range("F11").Value = "LogList@Group:V:0.2"
$ws.Range("C12").Value = "v,h;1-n"
$ws.Range("F12").Value = "Persons@Grid:10"

# --- Column F width (stored width becomes 22) ---
$ws.Columns("F").ColumnWidth = 21.25

# --- Selection moved to C13 ---
$ws.Range("C13").Select()
